$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.461.50'
$ws.Range('E2').Value = '  -2.82%  '
$ws.Range('D3').Value = '1.986.78'
$ws.Range('E3').Value = '  -3.38%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.19'
$ws.Range('E5').Value = '  -2.90%  '
$ws.Range('E6').Value = '  -3.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.06'
$ws.Range('E7').Value = '  -11.10%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.378'
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.66'
$ws.Range('E10').Value = '  -3.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0818'
$ws.Range('E11').Value = '  +6.49%  '
$ws.Range('E12').Value = '  -0.91%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.27'
$ws.Range('E13').Value = '  +14.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.870'
$ws.Range('E14').Value = '  -4.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.11'
$ws.Range('E15').Value = '  -5.84%  '
$ws.Range('D16').Value = '2.275.81'
$ws.Range('E16').Value = '  -3.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.48'
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('D18').Value = '1.991.30'
$ws.Range('E18').Value = '  -3.21%  '
$ws.Range('D19').Value = '36.421.65'
$ws.Range('E19').Value = '  -2.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.20'
$ws.Range('E20').Value = '  -3.63%  '
$ws.Range('D21').Value = '0.0₃0867'
$ws.Range('E21').Value = '  -1.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.35'
$ws.Range('E22').Value = '  -2.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.14'
$ws.Range('E23').Value = '  -2.16%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('E25').Value = '  -0.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.21'
$ws.Range('E27').Value = '  +4.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.34'
$ws.Range('E28').Value = '  +0.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.91'
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.125'
$ws.Range('E30').Value = '  +8.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.121'
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.20'
$ws.Range('E32').Value = '  -0.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.94'
$ws.Range('E33').Value = '  -6.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0635'
$ws.Range('E34').Value = '  +2.12%  '
$ws.Range('E35').Value = '  -6.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.33'
$ws.Range('E36').Value = '  +3.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.29'
$ws.Range('E37').Value = '  -6.33%  '
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('E39').Value = '  -3.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.11'
$ws.Range('E40').Value = '  +3.32%  '
$ws.Range('E41').Value = '  +1.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0971'
$ws.Range('E42').Value = '  -5.70%  '
$ws.Range('E43').Value = '  -3.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0215'
$ws.Range('E44').Value = '  -2.49%  '
$ws.Range('E45').Value = '  -4.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.42'
$ws.Range('E46').Value = '  -3.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '93.25'
$ws.Range('E47').Value = '  -2.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.65'
$ws.Range('E48').Value = '  -3.92%  '
$ws.Range('D49').Value = '1.375.96'
$ws.Range('E49').Value = '  -3.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.87'
$ws.Range('E50').Value = '  -2.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '45.41'
$ws.Range('E51').Value = '  -2.87%  '
